# 自动更新Excel文件 - 每日剩余天数递减脚本
# 逻辑：以"开始时间"(F列)为基准，"总天"(D列)减去从开始时间到当前日期经过的天数，
# 得到"剩余"(E列)。每天运行一次，相当于把每行的剩余天数减 1。
# 如果剩余天数已经耗尽(<=0)，代表该店铺续费，剩余天数重置为总天数，
# 并将开始时间刷新为当前参考日期。

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 当前参考日期（比上一次更新晚一天）
$newToday = 20251221

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row   # xlUp = -4162
if ($lastRow -lt 2) { $lastRow = 99 }

for ($r = 2; $r -le $lastRow; $r++) {
    $dCell = $ws.Cells.Item($r, 4)   # D列：总天
    $eCell = $ws.Cells.Item($r, 5)   # E列：剩余
    $fCell = $ws.Cells.Item($r, 6)   # F列：开始时间

    $dVal = $dCell.Value2
    $eVal = $eCell.Value2
    $fVal = $fCell.Value2

    if ($dVal -eq $null -or $eVal -eq $null -or $fVal -eq $null) {
        continue
    }

    # 跳过开始时间格式异常的行（无法解析的日期，数据本身有问题，保持不变）
    $fStr = [string]([int]$fVal)
    if ($fStr.Length -ne 8) {
        continue
    }

    $newE = [int]$eVal - 1

    if ($newE -le 0) {
        # 剩余天数耗尽 -> 续费：重置剩余为总天数，开始时间更新为当前参考日期
        $eCell.Value2 = [int]$dVal
        $fCell.Value2 = $newToday
    } else {
        $eCell.Value2 = $newE
    }
}
